$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column O (DÉLKA_PRACOVNÍHO_POMĚRU) first, then column M (ZAŘAZENO),
# so that indices of earlier columns remain stable while deleting.
$ws.Range("O1").EntireColumn.Delete()
$ws.Range("M1").EntireColumn.Delete()

# Update the view/selection to match the saved state in the diff.
$ws.Range("L13").Select()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
